$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ20846750",
    "summ20920454",
    "summ20998312",
    "summ21073729",
    "summ21156118",
    "summ21222585",
    "summ21298376",
    "summ21372596",
    "summ21453802",
    "summ21522902",
    "summ21598946",
    "summ21674803",
    "summ21743307",
    "summ21822833",
    "summ21898756",
    "summ21979117"
)

for ($i = 1; $i -le $newNames.Length; $i++) {
    $wb.Worksheets.Item($i).Name = $newNames[$i - 1]
}
